$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.479.50"
$ws.Range("E2").Value = "  +6.40%  "
$ws.Range("D3").Value = "2.654.13"
$ws.Range("E3").Value = "  +8.15%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "185.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "582.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.190"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.68%  "
$ws.Range("D10").Value = "2.650.87"
$ws.Range("E10").Value = "  +8.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +5.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "3.147.33"
$ws.Range("E14").Value = "  +8.08%  "
$ws.Range("D15").Value = "74.312.10"
$ws.Range("E15").Value = "  +6.32%  "
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("E17").Value = "  +8.79%  "
$ws.Range("D18").Value = "2.656.78"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +28.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.71%  "
$ws.Range("D29").Value = "2.792.77"
$ws.Range("E29").Value = "  +7.91%  "
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "0.0₃0931"
$ws.Range("E31").Value = "  +8.18%  "
$ws.Range("E32").Value = "  +12.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "517.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.43%  "
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("E35").Value = "  +6.24%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("E38").Value = "  +4.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "169.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +25.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.326"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.79%  "
$ws.Range("E49").Value = "  +15.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +20.63%  "
